$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.705.08'
$ws.Range('E2').Value = '  +1.36%  '
$ws.Range('D3').Value = '1.633.61'
$ws.Range('E3').Value = '  +0.72%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.60'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E6').Value = '  +3.17%  '
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('E8').Value = '  +1.25%  '
$ws.Range('E9').Value = '  +0.98%  '
$ws.Range('E10').Value = '  +2.39%  '
$ws.Range('E11').Value = '  +3.45%  '
$ws.Range('D12').Value = '1.860.70'
$ws.Range('E12').Value = '  +0.72%  '
$ws.Range('D13').Value = '1.654.54'
$ws.Range('E13').Value = '  +1.82%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.11'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.83%  '
$ws.Range('E15').Value = '  +0.81%  '
$ws.Range('D16').Value = '26.695.53'
$ws.Range('E16').Value = '  +1.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.58'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.69%  '
$ws.Range('E18').Value = '  +2.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '218.75'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +8.18%  '
$ws.Range('E20').Value = '  +0.19%  '
$ws.Range('E21').Value = '  +1.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.37'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.81%  '
$ws.Range('E24').Value = '  +5.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '148.00'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.45%  '
$ws.Range('E26').Value = '  +0.15%  '
$ws.Range('E27').Value = '  +1.39%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.87'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.13%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.55'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.35%  '
$ws.Range('E30').Value = '  -2.15%  '
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('E32').Value = '  +3.49%  '
$ws.Range('E33').Value = '  +2.39%  '
$ws.Range('E34').Value = '  +0.65%  '
$ws.Range('B35').Value = 'Maker'
$ws.Range('C35').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D35').Value = '1.226.00'
$ws.Range('E35').Value = '  +5.43%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.40'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.21%  '
$ws.Range('E37').Value = '  +5.74%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.807'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.68%  '
$ws.Range('E39').Value = '  +0.15%  '
$ws.Range('E40').Value = '  +1.17%  '
$ws.Range('E41').Value = '  -1.68%  '
$ws.Range('E42').Value = '  +1.72%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.36'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.68%  '
$ws.Range('D44').Value = '1.769.49'
$ws.Range('E44').Value = '  +0.57%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '92.66'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.54%  '
$ws.Range('E46').Value = '  +2.88%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.40'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.75%  '
$ws.Range('E48').Value = '  -0.28%  '
$ws.Range('E49').Value = '  +0.92%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.67'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.63%  '
$ws.Range('E51').Value = '  +0.01%  '
